$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (masthead volume/number + week-of dates) ---
$ws.Range("A8").Value = "Volume 31   Number  11"
$ws.Range("C9").Value = "Report Covering the Week  3/11/2024  Through  3/17/2024"

# --- Crime statistics table updates (rows 14-30) ---
# Row 14
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "0"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = "***.*"
$ws.Cells.Item(14,14).Value = -84.615384615384

# Row 15
$ws.Cells.Item(15,13).Value = -83.333333333333
$ws.Cells.Item(15,14).Value = -92.857142857142

# Row 16
$ws.Cells.Item(16,3).Value = 4
$ws.Cells.Item(16,4).Value = 1
$ws.Cells.Item(16,5).Value = 300
$ws.Cells.Item(16,6).Value = 12
$ws.Cells.Item(16,7).Value = 7
$ws.Cells.Item(16,8).Value = 71.428571428571
$ws.Cells.Item(16,9).Value = 29
$ws.Cells.Item(16,10).Value = 21
$ws.Cells.Item(16,11).Value = 38.095238095238
$ws.Cells.Item(16,12).Value = 45
$ws.Cells.Item(16,13).Value = -21.621621621621
$ws.Cells.Item(16,14).Value = -82.530120481927

# Row 17
$ws.Cells.Item(17,3).Value = 4
$ws.Cells.Item(17,4).Value = 5
$ws.Cells.Item(17,5).Value = -20
$ws.Cells.Item(17,6).Value = 14
$ws.Cells.Item(17,7).Value = 17
$ws.Cells.Item(17,8).Value = -17.647058823529
$ws.Cells.Item(17,9).Value = 37
$ws.Cells.Item(17,10).Value = 44
$ws.Cells.Item(17,11).Value = -15.909090909090
$ws.Cells.Item(17,12).Value = -21.276595744680
$ws.Cells.Item(17,13).Value = 54.166666666666
$ws.Cells.Item(17,14).Value = -71.317829457364

# Row 18
$ws.Cells.Item(18,3).Value = 3
$ws.Cells.Item(18,4).Value = 1
$ws.Cells.Item(18,5).Value = 200
$ws.Cells.Item(18,6).Value = 6
$ws.Cells.Item(18,7).Value = 7
$ws.Cells.Item(18,8).Value = -14.285714285714
$ws.Cells.Item(18,9).Value = 15
$ws.Cells.Item(18,10).Value = 18
$ws.Cells.Item(18,11).Value = -16.666666666666
$ws.Cells.Item(18,12).Value = -31.818181818181
$ws.Cells.Item(18,13).Value = 7.142857142857
$ws.Cells.Item(18,14).Value = -91.428571428571

# Row 19
$ws.Cells.Item(19,3).Value = 6
$ws.Cells.Item(19,4).Value = 6
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,6).Value = 19
$ws.Cells.Item(19,7).Value = 18
$ws.Cells.Item(19,8).Value = 5.555555555555
$ws.Cells.Item(19,9).Value = 73
$ws.Cells.Item(19,10).Value = 58
$ws.Cells.Item(19,11).Value = 25.862068965517
$ws.Cells.Item(19,12).Value = 5.797101449275
$ws.Cells.Item(19,13).Value = 180.769230769231
$ws.Cells.Item(19,14).Value = -9.876543209876

# Row 20
$ws.Cells.Item(20,3).Value = 1
$ws.Cells.Item(20,4).NumberFormat = '#,##0'
$ws.Cells.Item(20,4).Value = 2
$ws.Cells.Item(20,5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(20,5).Value = -50
$ws.Cells.Item(20,6).Value = 6
$ws.Cells.Item(20,7).Value = 5
$ws.Cells.Item(20,8).Value = 20
$ws.Cells.Item(20,9).Value = 17
$ws.Cells.Item(20,10).Value = 21
$ws.Cells.Item(20,11).Value = -19.047619047619
$ws.Cells.Item(20,12).Value = 70
$ws.Cells.Item(20,13).Value = 88.888888888888
$ws.Cells.Item(20,14).Value = -69.642857142857

# Row 21
$ws.Cells.Item(21,3).Value = 18
$ws.Cells.Item(21,4).Value = 15
$ws.Cells.Item(21,5).Value = 20
$ws.Cells.Item(21,6).Value = 57
$ws.Cells.Item(21,7).Value = 56
$ws.Cells.Item(21,8).Value = 1.785714285714
$ws.Cells.Item(21,9).Value = 174
$ws.Cells.Item(21,10).Value = 164
$ws.Cells.Item(21,11).Value = 6.097560975609
$ws.Cells.Item(21,12).Value = 2.352941176470
$ws.Cells.Item(21,13).Value = 48.717948717948
$ws.Cells.Item(21,14).Value = -72.555205047318

# Row 22
$ws.Cells.Item(22,3).NumberFormat = '#,##0'
$ws.Cells.Item(22,3).Value = 1
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "0"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = "***.*"
$ws.Cells.Item(22,6).NumberFormat = '#,##0'
$ws.Cells.Item(22,6).Value = 1
$ws.Cells.Item(22,7).Value = 3
$ws.Cells.Item(22,8).Value = -66.666666666666
$ws.Cells.Item(22,9).Value = 3
$ws.Cells.Item(22,10).Value = 10
$ws.Cells.Item(22,11).Value = -70
$ws.Cells.Item(22,12).Value = -40
$ws.Cells.Item(22,13).Value = 0

# Row 23
$ws.Cells.Item(23,10).Value = 4
$ws.Cells.Item(23,11).Value = -75

# Row 24
$ws.Cells.Item(24,3).Value = 20
$ws.Cells.Item(24,4).Value = 16
$ws.Cells.Item(24,5).Value = 25
$ws.Cells.Item(24,6).Value = 64
$ws.Cells.Item(24,7).Value = 57
$ws.Cells.Item(24,8).Value = 12.280701754386
$ws.Cells.Item(24,9).Value = 164
$ws.Cells.Item(24,10).Value = 158
$ws.Cells.Item(24,11).Value = 3.797468354430
$ws.Cells.Item(24,12).Value = -9.392265193370
$ws.Cells.Item(24,13).Value = 192.857142857143

# Row 25
$ws.Cells.Item(25,3).Value = 5
$ws.Cells.Item(25,4).Value = 1
$ws.Cells.Item(25,5).Value = 400
$ws.Cells.Item(25,6).Value = 13
$ws.Cells.Item(25,7).Value = 10
$ws.Cells.Item(25,8).Value = 30
$ws.Cells.Item(25,9).Value = 31
$ws.Cells.Item(25,10).Value = 40
$ws.Cells.Item(25,11).Value = -22.5
$ws.Cells.Item(25,12).Value = -56.338028169014

# Row 26
$ws.Cells.Item(26,3).Value = 5
$ws.Cells.Item(26,4).Value = 5
$ws.Cells.Item(26,5).Value = 0
$ws.Cells.Item(26,6).Value = 28
$ws.Cells.Item(26,8).Value = 27.272727272727
$ws.Cells.Item(26,9).Value = 62
$ws.Cells.Item(26,10).Value = 58
$ws.Cells.Item(26,11).Value = 6.896551724137
$ws.Cells.Item(26,12).Value = -24.390243902439
$ws.Cells.Item(26,13).Value = -30.337078651685

# Row 27
$ws.Cells.Item(27,12).Value = -25

# Row 28
$ws.Cells.Item(28,12).Value = -53.846153846153

# Row 29
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "0"
$ws.Cells.Item(29,5).NumberFormat = "@"
$ws.Cells.Item(29,5).Value = "***.*"
$ws.Cells.Item(29,14).Value = -95.652173913043

# Row 30
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "0"
$ws.Cells.Item(30,5).NumberFormat = "@"
$ws.Cells.Item(30,5).Value = "***.*"
$ws.Cells.Item(30,14).Value = -95.652173913043
